$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 0
$ws1.Range("F3").Value = 0
$ws1.Range("F4").Value = 4853
$ws1.Range("F5").Value = 212
$ws1.Range("F6").Value = 162
$ws1.Range("F8").Value = 113
$ws1.Range("F9").Value = 98
$ws1.Range("F12").Value = 0
$ws1.Range("F13").Value = 0
$ws1.Range("F14").Value = 258
$ws1.Range("F15").Value = 0
$ws1.Range("F18").Value = 158
$ws1.Range("F19").Value = 0
$ws1.Range("F21").Value = 0
$ws1.Range("F23").Value = 41
$ws1.Range("F25").Value = 0
$ws1.Range("F27").Value = 4005
$ws1.Range("F28").Value = 411
$ws1.Range("F30").Value = 0
$ws1.Range("F31").Value = 0
$ws1.Range("F33").Value = 536
$ws1.Range("F34").Value = 152
$ws1.Range("F35").Value = 305
$ws1.Range("F38").Value = 185
$ws1.Range("F39").Value = 13
$ws1.Range("F43").Value = 81
$ws1.Range("F45").Value = 503
$ws1.Range("F46").Value = 0
$ws1.Range("F47").Value = 6
$ws1.Range("F49").Value = 593

# Sheet 2: 演出
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 112

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 37
$ws4.Range("F3").Value = 0
$ws4.Range("F4").Value = 4853
$ws4.Range("F5").Value = 0
$ws4.Range("F8").Value = 112
$ws4.Range("F10").Value = 98
$ws4.Range("F11").Value = 766
$ws4.Range("F12").Value = 0
$ws4.Range("F13").Value = 0
$ws4.Range("F14").Value = 0
$ws4.Range("F15").Value = 0
$ws4.Range("F18").Value = 158
$ws4.Range("F21").Value = 0
$ws4.Range("F23").Value = 41
$ws4.Range("F24").Value = 88
$ws4.Range("F25").Value = 545
$ws4.Range("F27").Value = 0
$ws4.Range("F29").Value = 52
$ws4.Range("F30").Value = 0
$ws4.Range("F32").Value = 569
$ws4.Range("F33").Value = 536
$ws4.Range("F34").Value = 152
$ws4.Range("F35").Value = 305
$ws4.Range("F36").Value = 0
$ws4.Range("F37").Value = 381
$ws4.Range("F40").Value = 1576
$ws4.Range("F41").Value = 980
$ws4.Range("F42").Value = 48
$ws4.Range("F43").Value = 81
$ws4.Range("F46").Value = 484
$ws4.Range("F49").Value = 593
